$wb = $excel.ActiveWorkbook

# Rename sheets (workbook.xml sheet name changes)
$wb.Worksheets.Item(1).Name = "GNG_TO-16511686534598207"
$wb.Worksheets.Item(2).Name = "NB_TO-16511686573367705"
$wb.Worksheets.Item(3).Name = "RS_TO-1651168657338707"
$wb.Worksheets.Item(4).Name = "TOL_TO-16511686573861086"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16511686574622436"

# Sheet 1 (GNG)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16511686534218023.csv"
$ws1.Range("B3").Value = "GNG_stims-16511686534426196.csv"
$ws1.Range("B4").Value = "go_stims-16511686534456115.csv"
$ws1.Range("B5").Value = "GNG_stims-1651168653458823.csv"

# Sheet 2 (NB)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-16511686547535627.csv"
$ws2.Range("B3").Value = "TB-16511686573114028.csv"
$ws2.Range("B4").Value = "ZB-match_9-16511686538366807.csv"
$ws2.Range("B5").Value = "TB-16511686566443732.csv"
$ws2.Range("B6").Value = "OB-16511686542225082.csv"
$ws2.Range("B7").Value = "ZB-match_2-16511686535152845.csv"
$ws2.Range("B8").Value = "OB-16511686543955374.csv"
$ws2.Range("B9").Value = "TB-1651168656971683.csv"
$ws2.Range("B10").Value = "ZB-match_2-16511686537483432.csv"

# Sheet 3 (RS)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes open"
$ws3.Range("B3").Value = "eyes closed"

# Sheet 4 (TOL)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16511686573525264.csv"
$ws4.Range("B3").Value = "ZM_stims-16511686573407562.csv"
$ws4.Range("B4").Value = "MM_stims-1651168657367715.csv"
$ws4.Range("B5").Value = "ZM_stims-1651168657353494.csv"
$ws4.Range("B6").Value = "MM_stims-16511686573841095.csv"
$ws4.Range("B7").Value = "ZM_stims-16511686573687053.csv"

# Sheet 5 (vSAT)
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-1651168657416175.csv"
$ws5.Range("B3").Value = "vSAT_stims-16511686574469657.csv"
$ws5.Range("B4").Value = "SAT_stims-16511686573928082.csv"
$ws5.Range("B5").Value = "vSAT_stims-16511686574308195.csv"
